{"js": "// Remove <w:contextualSpacing .../> paragraph-property elements that were\n// left over from the previous authoring tool: every w:pPr in this document\n// (and in word/comments.xml) currently carries a redundant\n// <w:contextualSpacing w:val=\"0\"/>, which the commit removes everywhere.\n//\n// The Word JS API does not expose a `contextualSpacing` member on\n// ParagraphFormat, so we round-trip the package OOXML: pull the full\n// flat-OPC package, strip every <w:contextualSpacing/> element via a\n// surgical string replace (this does not touch any other markup), and\n// write the package back with Range.insertOoxml (Replace).\nconst body = context.document.body;\nconst ooxml = body.getOoxml();\nawait context.sync();\n\nconst original = ooxml.value;\nconst updated = original.replace(/<w:contextualSpacing\\b[^>]*\\/>/g, \"\");\n\nif (updated !== original) {\n  body.insertOoxml(updated, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Remove <w:contextualSpacing .../> paragraph-property elements that were\n# left over from the previous authoring tool: every w:pPr in this document\n# (and in word/comments.xml) currently carries a redundant\n# <w:contextualSpacing w:val=\"0\"/>, which this change removes everywhere.\n#\n# The Word object model does not expose a ContextualSpacing member on\n# ParagraphFormat/Paragraph, so we round-trip the package OOXML: read the\n# full flat-OPC package from Range.WordOpenXML, strip every\n# <w:contextualSpacing/> element via a surgical string replace (this does\n# not touch any other markup), and write the package back with\n# Range.InsertXML.\n$d = $word.ActiveDocument\n$rng = $d.Content\n$xml = $rng.WordOpenXML\n\n$pattern = '<w:contextualSpacing\\b[^>]*/>'\n$updated = [System.Text.RegularExpressions.Regex]::Replace($xml, $pattern, '')\n\nif ($updated -ne $xml) {\n  $null = $rng.InsertXML($updated)\n}\n"}
